$d = $word.ActiveDocument

# --- 1. Split "Aangepast: posities..." run after "Aan" and drop a _GoBack bookmark there ---
$findRange = $d.Content
$findRange.Find.Execute("Aangepast: posities worden weggelaten")
$splitPos = $findRange.Start + 3
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 2. Append a new row to the (only) table with the B010 bug entry ---
$t = $d.Tables(1)
$newRow = $t.Rows.Add()
$cells = $newRow.Cells

$cells.Item(1).Range.Text = "13 dec 2020"
$cells.Item(2).Range.Text = "B010"
$cells.Item(3).Range.Text = "Verkeerde portefeuille opgeslagen"
$cells.Item(4).Range.Text = "Als je portefeuille01.csv probeert op te slaan wordt portefeuille.csv daarmee overschreven"
$cells.Item(5).Range.Text = "Hoog"

# Cell 6 ("Open 13 dec") needs two runs sharing identical (theme) formatting
$c6 = $cells.Item(6)
$c6.Range.Text = "Open"
$c6start = $c6.Range.Start
$c6end = $c6.Range.End
$insertPoint = $d.Range($c6end - 1, $c6end - 1)
$insertPoint.InsertAfter(" 13 dec")

$wholeC6 = $d.Range($c6start, $c6start + 11)
$wholeC6.Font.TextColor.ObjectThemeColor = 13

# force the run split between "Open" and " 13 dec" (both ends up with identical rPr)
$secondPart = $d.Range($c6start + 4, $c6start + 11)
$secondPart.Bold = 1
$secondPart.Bold = 0

$cells.Item(7).Range.Text = "Voorlopige voorziening om meerdere portefeuilles te gebruiken (nl in IntelliJ aan component toevoegen en handmatig kopieren naar bestand) is zo niet mogelijk"
$cells.Item(8).Range.Text = "037"
